$wb = $excel.ActiveWorkbook

# The "Logical Operators" sheet gets the new formulas in D2:G16
$ws = $wb.Worksheets.Item("Logical Operators")

# Row 2 - plain (non-shared) formulas
$ws.Range("D2").Formula = '=IF(C2>=90, "PASS", "FAIL")'
$ws.Range("E2").Formula = '=IF(C2>=90,"A",IF(AND(C2>=80,C2<=89),"B",IF(AND(C2>=70,C2<=79),"C",IF(AND(C2>=60,C2<=69),"D","F"))))'
$ws.Range("F2").Formula = '=IF(OR(C2<60, C2>90), "OUTLIER", "AVG")'
$ws.Range("G2").Formula = '=IF(AND(C2>95, B2="M"), "Male Achiever", IF(AND(C2>95, B2="F"), "Female Achiever", "None"))'

# Rows 3-16 - shared formulas (filled down from row 3)
$ws.Range("D3:D16").Formula = '=IF(C3>=90, "PASS", "FAIL")'
$ws.Range("E3:E16").Formula = '=IF(C3>=90,"A",IF(AND(C3>=80,C3<=89),"B",IF(AND(C3>=70,C3<=79),"C",IF(AND(C3>=60,C3<=69),"D","F"))))'
$ws.Range("F3:F16").Formula = '=IF(OR(C3<60, C3>90), "OUTLIER", "AVG")'
$ws.Range("G3:G16").Formula = '=IF(AND(C3>95, B3="M"), "Male Achiever", IF(AND(C3>95, B3="F"), "Female Achiever", "None"))'

# Update the active sheet / tab selection: "Logical Operators" becomes the
# active, selected tab (with G4 selected), while "Formulas 101" loses its
# tabSelected flag.
$ws.Range("G4").Select()
$ws.Activate()

$wb.Save()
